# Apply cell value updates per the diff (odds/stat columns for several match rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6
$ws.Range("G6").Value = 3.3
$ws.Range("I6").Value = 2.5
$ws.Range("J6").Value = 1.14
$ws.Range("K6").Value = 5.5
$ws.Range("N6").Value = 3.1
$ws.Range("O6").Value = 1.36
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 2.38
$ws.Range("S6").Value = 1.53
$ws.Range("V6").Value = 13
$ws.Range("AB6").Value = 21
$ws.Range("AE6").Value = 10
$ws.Range("AG6").Value = 26

# Row 10
$ws.Range("G10").Value = 2.38
$ws.Range("I10").Value = 2.7
$ws.Range("J10").Value = 1.04
$ws.Range("K10").Value = 12
$ws.Range("T10").Value = 9
$ws.Range("U10").Value = 12
$ws.Range("V10").Value = 9.5
$ws.Range("W10").Value = 23
$ws.Range("X10").Value = 19
$ws.Range("Y10").Value = 26
$ws.Range("Z10").Value = 12
$ws.Range("AA10").Value = 7
$ws.Range("AE10").Value = 15
$ws.Range("AG10").Value = 29

# Row 12
$ws.Range("G12").Value = 1.91
$ws.Range("J12").Value = 1.1
$ws.Range("K12").Value = 7
$ws.Range("L12").Value = 1.44
$ws.Range("M12").Value = 2.63
$ws.Range("N12").Value = 2.35
$ws.Range("O12").Value = 1.57
$ws.Range("P12").Value = 1.53
$ws.Range("Q12").Value = 2.38
$ws.Range("R12").Value = 2.1
$ws.Range("S12").Value = 1.67
$ws.Range("X12").Value = 19
$ws.Range("Z12").Value = 7
$ws.Range("AA12").Value = 6
$ws.Range("AB12").Value = 19
$ws.Range("AG12").Value = 41
$ws.Range("AI12").Value = 51

# Row 13
$ws.Range("L13").Value = 1.36
$ws.Range("M13").Value = 3

# Row 15
$ws.Range("G15").Value = 1.6
$ws.Range("H15").Value = 3.5
$ws.Range("I15").Value = 6.5
$ws.Range("N15").Value = 2.5
$ws.Range("O15").Value = 1.5
$ws.Range("T15").Value = 4.75
$ws.Range("U15").Value = 6
$ws.Range("W15").Value = 11
$ws.Range("X15").Value = 17
$ws.Range("AA15").Value = 7.5
$ws.Range("AB15").Value = 26
$ws.Range("AE15").Value = 29
$ws.Range("AG15").Value = 81

# Row 16
$ws.Range("G16").Value = 1.95
$ws.Range("H16").Value = 3.25
$ws.Range("I16").Value = 4.1
$ws.Range("N16").Value = 2.4
$ws.Range("O16").Value = 1.53
$ws.Range("V16").Value = 9.5
$ws.Range("W16").Value = 17
$ws.Range("AD16").Value = 9
$ws.Range("AE16").Value = 19

# Row 18
$ws.Range("N18").Value = 2
$ws.Range("O18").Value = 1.8

# Row 44
$ws.Range("G44").Value = 3.8
$ws.Range("H44").Value = 3.5
$ws.Range("I44").Value = 1.91
$ws.Range("L44").Value = 1.29
$ws.Range("M44").Value = 3.5
$ws.Range("N44").Value = 1.95
$ws.Range("O44").Value = 1.9
$ws.Range("X44").Value = 29
$ws.Range("AA44").Value = 6.5
$ws.Range("AC44").Value = 41
$ws.Range("AE44").Value = 9.5
$ws.Range("AG44").Value = 17

# Row 45
$ws.Range("G45").Value = 2.62
$ws.Range("H45").Value = 3.25
$ws.Range("I45").Value = 2.5
$ws.Range("M45").Value = 5
$ws.Range("R45").Value = 1.38
$ws.Range("S45").Value = 2.6
$ws.Range("T45").Value = 14.5
$ws.Range("U45").Value = 19
$ws.Range("W45").Value = 35
$ws.Range("X45").Value = 18.5
$ws.Range("Y45").Value = 18.5
$ws.Range("Z45").Value = 15
$ws.Range("AA45").Value = 7
$ws.Range("AB45").Value = 9.75
$ws.Range("AC45").Value = 29
$ws.Range("AD45").Value = 12
$ws.Range("AG45").Value = 30
$ws.Range("AH45").Value = 18

# Row 49
$ws.Range("L49").Value = 1.29
$ws.Range("M49").Value = 3.5
$ws.Range("N49").Value = 2
$ws.Range("O49").Value = 1.85
